{"js": "// Fix the typo \"on Internet\" -> \"on the Internet\" in the walkthrough text.\nconst searchResults = context.document.body.search(\"on Internet and I found a possible exploit\", { matchCase: true });\nsearchResults.load(\"items,text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Target text \"on Internet and I found a possible exploit\" not found.');\n}\n\nconst range = searchResults.items[0];\nrange.insertText(\"on the Internet and I found a possible exploit\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Correct a mistake: \"on Internet\" into \"on the Internet\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"on Internet and I found a possible exploit\"\n$find.Replacement.Text = \"on the Internet and I found a possible exploit\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
